# Upload new version with timestamp
# - Adds a new stock-shortage line item "UNITREXATE 50 MG 5 I.M. I.V. VIALS"
#   as a new row 30 (pushing the existing rows 30-39 down to 31-40).
# - Updates the running total accordingly.
# - Updates the "printed at" timestamp footer string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row at position 30 (shifts rows 30.. down by one, carrying
#    along the totals row and the footer row automatically).
# ---------------------------------------------------------------------------
$ws.Rows("30:30").Insert()

# Copy the (now shifted-down) row 31 formatting into the freshly inserted,
# still-blank row 30 so the new product row looks like every other product
# row in the table (borders / fills / number formats / fonts).
$ws.Range("A31:Q31").Copy()
$ws.Range("A30:Q30").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Re-create the merged cells that make up one product row.
$ws.Range("A30:B30").Merge()
$ws.Range("C30:G30").Merge()
$ws.Range("H30:K30").Merge()
$ws.Range("L30:M30").Merge()
$ws.Range("N30:O30").Merge()

# Match the row height used by the product rows around it.
$ws.Rows("30:30").RowHeight = 24.75

# ---------------------------------------------------------------------------
# 2. Fill in the new item's data.
# ---------------------------------------------------------------------------
$ws.Range("A30").Value = 24
$ws.Range("C30").Value = "UNITREXATE 50 MG 5 I.M. I.V. VIALS"
$ws.Range("H30").Value = "1:0"
$ws.Range("L30").Value = "1"
$ws.Range("N30").Value = "385.00"
$ws.Range("P30").Value = "77.0000"
$ws.Range("Q30").Value = "0:1"

# ---------------------------------------------------------------------------
# 3. Update the grand-total row (now row 39 after the insert) to include the
#    new item's sale price (1118.98 + 77.00 = 1195.98).
# ---------------------------------------------------------------------------
$ws.Range("P39").Value = 1195.98

# ---------------------------------------------------------------------------
# 4. Update the "printed at" timestamp shown in the footer.
# ---------------------------------------------------------------------------
$ws.Range("A40").Value = "Thursday, 7 August, 2025 12:48 PM"

Write-Output "edit complete"
